# Renouveau DonneeTest et Changement commentaire
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values (columns A, B, C, E, F, G, H) - I2 stays 3 (unchanged)
$ws.Range("A2").Value = 70
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 65
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 35
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 37

# New weight values for column C, rows 3 through 71
$cValues = @(26,44,93,68,31,54,53,91,32,9,31,19,74,8,21,35,68,14,23,6,5,22,70,54,84,2,95,74,63,26,28,97,45,6,19,3,79,11,91,65,44,99,48,87,89,10,31,17,87,80,98,15,35,9,61,96,93,21,88,56,30,18,94,89,34,79,10,6,97)

$row = 3
foreach ($val in $cValues) {
    $ws.Cells.Item($row, 3).Value = $val
    $row++
}
